$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and two row re-sorts) per the
# Fri Apr 28 18:56:11 UTC 2023 GitHub Actions refresh of cryptos list.

# Row 2
$ws.Range("D2").Value = "29.246.42"
$ws.Range("E2").Value = "  -1.63%  "

# Row 3
$ws.Range("D3").Value = "1.893.61"
$ws.Range("E3").Value = "  -1.75%  "

# Row 4
$ws.Range("E4").Value = "  +0.34%  "

# Row 5
$ws.Range("D5").Value = "323.22"
$ws.Range("E5").Value = "  -3.43%  "

# Row 6
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.34%  "

# Row 7
$ws.Range("D7").Value = "0.4761"
$ws.Range("E7").Value = "  +2.00%  "

# Row 8
$ws.Range("D8").Value = "0.4047"
$ws.Range("E8").Value = "  -2.70%  "

# Row 9
$ws.Range("D9").Value = "0.08019"
$ws.Range("E9").Value = "  -0.44%  "

# Row 10
$ws.Range("D10").Value = "0.9982"
$ws.Range("E10").Value = "  -2.56%  "

# Row 11
$ws.Range("D11").Value = "23.20"
$ws.Range("E11").Value = "  +3.87%  "

# Row 12
$ws.Range("D12").Value = "1.976.84"
$ws.Range("E12").Value = "  +2.34%  "

# Row 13
$ws.Range("D13").Value = "5.921"
$ws.Range("E13").Value = "  -1.49%  "

# Row 14
$ws.Range("D14").Value = "7.026"
$ws.Range("E14").Value = "  -2.25%  "

# Row 15
$ws.Range("D15").Value = "89.18"
$ws.Range("E15").Value = "  -0.79%  "

# Row 16
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.33%  "

# Row 17
$ws.Range("D17").Value = "0.06640"
$ws.Range("E17").Value = "  +0.71%  "

# Row 18
$ws.Range("E18").Value = "  -0.74%  "

# Row 19
$ws.Range("D19").Value = "17.50"
$ws.Range("E19").Value = "  -1.85%  "

# Row 20
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.49%  "

# Row 21
$ws.Range("D21").Value = "29.247.64"
$ws.Range("E21").Value = "  -1.45%  "

# Row 22
$ws.Range("E22").Value = "  -0.67%  "

# Row 23
$ws.Range("D23").Value = "11.68"
$ws.Range("E23").Value = "  +1.23%  "

# Row 24
$ws.Range("D24").Value = "2.157"
$ws.Range("E24").Value = "  -2.05%  "

# Row 25
$ws.Range("D25").Value = "2.079.07"
$ws.Range("E25").Value = "  -3.65%  "

# Row 26
$ws.Range("D26").Value = "154.32"
$ws.Range("E26").Value = "  -1.88%  "

# Row 27
$ws.Range("D27").Value = "19.74"
$ws.Range("E27").Value = "  -1.05%  "

# Row 28
$ws.Range("D28").Value = "5.903"
$ws.Range("E28").Value = "  +4.07%  "

# Row 29
$ws.Range("D29").Value = "2.085"
$ws.Range("E29").Value = "  -3.82%  "

# Row 30
$ws.Range("D30").Value = "117.83"
$ws.Range("E30").Value = "  -0.12%  "

# Row 31
$ws.Range("D31").Value = "1.024"
$ws.Range("E31").Value = "  -1.43%  "

# Row 32
$ws.Range("D32").Value = "0.09409"
$ws.Range("E32").Value = "  -0.57%  "

# Row 33
$ws.Range("D33").Value = "3.527"
$ws.Range("E33").Value = "  -0.35%  "

# Row 34
$ws.Range("D34").Value = "1.375"
$ws.Range("E34").Value = "  -4.41%  "

# Row 35
$ws.Range("D35").Value = "5.341"
$ws.Range("E35").Value = "  -1.95%  "

# Row 36
$ws.Range("E36").Value = "  -1.09%  "

# Row 37
$ws.Range("D37").Value = "0.06029"
$ws.Range("E37").Value = "  -2.10%  "

# Row 38
$ws.Range("D38").Value = "1.168"
$ws.Range("E38").Value = "  -1.02%  "

# Row 39
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "0.5834"
$ws.Range("E39").Value = "  -1.54%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "7.848"
$ws.Range("E40").Value = "  -7.38%  "

# Row 41
$ws.Range("D41").Value = "0.1831"
$ws.Range("E41").Value = "  -0.56%  "

# Row 42
$ws.Range("D42").Value = "10.04"
$ws.Range("E42").Value = "  -2.30%  "

# Row 43
$ws.Range("E43").Value = "  +2.02%  "

# Row 44
$ws.Range("D44").Value = "0.07693"
$ws.Range("E44").Value = "  +1.84%  "

# Row 45
$ws.Range("D45").Value = "2.355"
$ws.Range("E45").Value = "  -0.62%  "

# Row 46
$ws.Range("D46").Value = "12.19"
$ws.Range("E46").Value = "  -0.23%  "

# Row 47
$ws.Range("D47").Value = "0.5479"
$ws.Range("E47").Value = "  -2.07%  "

# Row 48
$ws.Range("D48").Value = "1.909"
$ws.Range("E48").Value = "  -1.76%  "

# Row 49
$ws.Range("D49").Value = "112.96"
$ws.Range("E49").Value = "  +0.12%  "

# Row 50
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "43.57"
$ws.Range("E50").Value = "  -1.04%  "

# Row 51
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "0.2920"
$ws.Range("E51").Value = "  -2.28%  "
